$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('ECs', 'Ntn3', 'Cdon', 'ECs', 3.0, 1.0, 0.248739, 0.7462169999999999, 0.09500454514204523, 0.09500454514204522, 3.0, 1.0, 1.922246666666667, 5.76674, 0.04282684507083529, 0.04282684507083529, 0.47813771362, 4.30323942258, 0.004068744935823549, 0.004068744935823548),
    @('ECs', 'Ntn3', 'Cdon', 'FAPs', 3.0, 1.0, 0.248739, 0.7462169999999999, 0.09500454514204523, 0.09500454514204522, 3.0, 1.0, 21.84955866666667, 65.548676, 0.4867989525538483, 0.4867989525538483, 5.434837373188, 48.91353635869199, 0.04624811306300242, 0.04624811306300241),
    @('ECs', 'Ntn3', 'Cdon', 'M2', 3.0, 1.0, 0.248739, 0.7462169999999999, 0.09500454514204523, 0.09500454514204522, 3.0, 1.0, 0.102188, 0.306564, 0.002276705544605019, 0.002276705544605019, 0.025418140932, 0.228763268388, 0.0002162973746875722, 0.0002162973746875722),
    @('ECs', 'Ntn3', 'Cdon', 'sCs', 3.0, 1.0, 0.248739, 0.7462169999999999, 0.09500454514204523, 0.09500454514204522, 3.0, 1.0, 21.01015966666667, 63.030479, 0.4680974968307114, 0.4680974968307114, 5.226046105326999, 47.03441494794299, 0.0444713897685317, 0.0444713897685317),
    @('FAPs', 'Ntn3', 'Cdon', 'ECs', 3.0, 1.0, 1.982221333333333, 5.946664, 0.7570989516890869, 0.7570989516890868, 3.0, 1.0, 1.922246666666667, 5.76674, 0.04282684507083529, 0.04282684507083529, 3.810318350595556, 34.29286515536, 0.03242415950728034, 0.03242415950728034),
    @('FAPs', 'Ntn3', 'Cdon', 'FAPs', 3.0, 1.0, 1.982221333333333, 5.946664, 0.7570989516890869, 0.7570989516890868, 3.0, 1.0, 21.84955866666667, 65.548676, 0.4867989525538483, 0.4867989525538483, 43.31066131298489, 389.795951816864, 0.3685549766618641, 0.368554976661864),
    @('FAPs', 'Ntn3', 'Cdon', 'M2', 3.0, 1.0, 1.982221333333333, 5.946664, 0.7570989516890869, 0.7570989516890868, 3.0, 1.0, 0.102188, 0.306564, 0.002276705544605019, 0.002276705544605019, 0.2025592336106667, 1.823033102496, 0.001723691381125192, 0.001723691381125191),
    @('FAPs', 'Ntn3', 'Cdon', 'sCs', 3.0, 1.0, 1.982221333333333, 5.946664, 0.7570989516890869, 0.7570989516890868, 3.0, 1.0, 21.01015966666667, 63.030479, 0.4680974968307114, 0.4680974968307114, 41.64678670800622, 374.821080372056, 0.3543961241388173, 0.3543961241388173),
    @('M2', 'Ntn3', 'Cdon', 'ECs', 1.0, 0.3333333333333333, 0.07456866666666666, 0.223706, 0.02848110774151, 0.02848110774151, 3.0, 1.0, 1.922246666666667, 5.76674, 0.04282684507083529, 0.04282684507083529, 0.1433393709377778, 1.29005433844, 0.001219755988691417, 0.001219755988691417),
    @('M2', 'Ntn3', 'Cdon', 'FAPs', 1.0, 0.3333333333333333, 0.07456866666666666, 0.223706, 0.02848110774151, 0.02848110774151, 3.0, 1.0, 21.84955866666667, 65.548676, 0.4867989525538483, 0.4867989525538483, 1.629292457028444, 14.663632113256, 0.01386457341614037, 0.01386457341614037),
    @('M2', 'Ntn3', 'Cdon', 'M2', 1.0, 0.3333333333333333, 0.07456866666666666, 0.223706, 0.02848110774151, 0.02848110774151, 3.0, 1.0, 0.102188, 0.306564, 0.002276705544605019, 0.002276705544605019, 0.007620022909333332, 0.068580206184, 0.00006484309591158876, 0.00006484309591158876),
    @('M2', 'Ntn3', 'Cdon', 'sCs', 1.0, 0.3333333333333333, 0.07456866666666666, 0.223706, 0.02848110774151, 0.02848110774151, 3.0, 1.0, 21.01015966666667, 63.030479, 0.4680974968307114, 0.4680974968307114, 1.566699592797111, 14.100296335174, 0.01333193524076663, 0.01333193524076663),
    @('sCs', 'Ntn3', 'Cdon', 'ECs', 3.0, 1.0, 0.312651, 0.9379529999999999, 0.1194153954273579, 0.1194153954273579, 3.0, 1.0, 1.922246666666667, 5.76674, 0.04282684507083529, 0.04282684507083529, 0.60099234258, 5.40893108322, 0.005114184639039991, 0.005114184639039991),
    @('sCs', 'Ntn3', 'Cdon', 'FAPs', 3.0, 1.0, 0.312651, 0.9379529999999999, 0.1194153954273579, 0.1194153954273579, 3.0, 1.0, 21.84955866666667, 65.548676, 0.4867989525538483, 0.4867989525538483, 6.831286366691999, 61.481577300228, 0.05813128941284145, 0.05813128941284144),
    @('sCs', 'Ntn3', 'Cdon', 'M2', 3.0, 1.0, 0.312651, 0.9379529999999999, 0.1194153954273579, 0.1194153954273579, 3.0, 1.0, 0.102188, 0.306564, 0.002276705544605019, 0.002276705544605019, 0.03194918038799999, 0.287542623492, 0.0002718736928806666, 0.0002718736928806666),
    @('sCs', 'Ntn3', 'Cdon', 'sCs', 3.0, 1.0, 0.312651, 0.9379529999999999, 0.1194153954273579, 0.1194153954273579, 3.0, 1.0, 21.01015966666667, 63.030479, 0.4680974968307114, 0.4680974968307114, 6.568847429942999, 59.11962686948699, 0.05589804768259583, 0.05589804768259583)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
